$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("edit-fields")

# --- Insert the two new "menu" rows above the old row 15 (activity/item_id) ---
$ws.Rows.Item(15).EntireRow.Insert()
$ws.Rows.Item(15).EntireRow.Insert()

$ws.Range("A15").Value = "menu"
$ws.Range("B15").Value = "order"
$ws.Range("C15").Value = "display_list"
$ws.Range("D15").Value = "none"

$ws.Range("A16").Value = "menu"
$ws.Range("B16").Value = "name"
$ws.Range("C16").Value = "display_list"
$ws.Range("D16").Value = "show"

# --- Insert the new "variable" row after inbox/to_user_id (now row 20) ---
$ws.Rows.Item(21).EntireRow.Insert()

$ws.Range("A21").Value = "variable"
$ws.Range("B21").Value = "value"
$ws.Range("C21").Value = "display_list"
$ws.Range("D21").Value = "show"

# --- Match formatting used by the neighbouring rows for the new entries ---
$ws.Range("A15:B16").HorizontalAlignment = -4131
$ws.Range("A15:B16").VerticalAlignment = -4160
$ws.Range("A21:B21").HorizontalAlignment = -4131
$ws.Range("A21:B21").VerticalAlignment = -4160

# --- Selection shown in the saved workbook ---
$ws.Range("C21:D21").Select()
